# Add "2022-Q3" data: new worksheet inserted before "2022-Q2", plus a
# corresponding new row in the "总计" (summary) sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" sheet right before the existing "2022-Q2"
#    sheet (keeps overall sheet ordering: 总计, 2022-Q3, 2022-Q2, 2022-Q1, ...)
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

# Header row
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Fund rows for 2022-Q3. Columns B,C,D,E,F,G are stored as text (to match
# the source workbook convention), column A is a numeric 0-based index and
# H is a numeric rank.
$q3rows = @(
    @("005583", "易方达港股通红利灵活配置混合", "6.98", "90.31", "4.34", "0.3029", 5),
    @("009778", "长信消费升级混合A",           "1.99", "89.28", "4.70", "0.0935", 9),
    @("009779", "长信消费升级混合C",           "0.89", "89.28", "4.70", "0.0418", 9),
    @("001715", "工银新焦点灵活配置混合A",      "0.41", "83.76", "8.96", "0.0367", 5),
    @("519959", "长信多利灵活配置混合A",        "0.95", "89.23", "3.81", "0.0362", 9),
    @("001998", "工银新焦点灵活配置混合C",      "0.21", "83.76", "8.96", "0.0188", 5),
    @("015774", "长信多利灵活配置混合E",        "0.07", "89.23", "3.81", "0.0027", 9),
    @("013488", "长信多利灵活配置混合C",        "0.05", "89.23", "3.81", "0.0019", 9)
)

# Force text storage for the numeric-looking columns (B and D..G) so they
# don't get auto-coerced to numbers, matching the source file's use of
# inline/text strings for these columns.
$q3.Range("B2:B9").NumberFormat = "@"
$q3.Range("D2:G9").NumberFormat = "@"

for ($i = 0; $i -lt $q3rows.Count; $i++) {
    $r = 2 + $i
    $row = $q3rows[$i]
    $q3.Range("A$r").Value = $i
    $q3.Range("B$r").Value = $row[0]
    $q3.Range("C$r").Value = $row[1]
    $q3.Range("D$r").Value = $row[2]
    $q3.Range("E$r").Value = $row[3]
    $q3.Range("F$r").Value = $row[4]
    $q3.Range("G$r").Value = $row[5]
    $q3.Range("H$r").Value = $row[6]
}

# Drop the temporary text NumberFormat now that the values are locked in as
# text - ClearFormats keeps the stored type as text while removing the
# leftover style index, matching cells elsewhere in the workbook that carry
# no explicit style.
$q3.Range("B2:H9").ClearFormats()

# Replicate the bold/centered/bordered style used for the header row and
# column-A index cells elsewhere in the workbook (style index "2"), by
# copying formatting only from the equivalent cells on the summary sheet.
$summary.Range("B1:D1").Copy() | Out-Null
$q3.Range("B1:H1").PasteSpecial(-4122) | Out-Null
$summary.Range("A2").Copy() | Out-Null
$q3.Range("A2:A9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Insert a new row 2 in "总计" for the 2022-Q3 summary figures, and
#    renumber the existing index column so it stays a 0-based sequence.
# ---------------------------------------------------------------------
$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 8
$summary.Range("D2").Value = 0.53

# Re-apply the column-A style (Insert() leaves the new row's A-cell
# unstyled) and renumber A3:A8 as 1..6.
$summary.Range("A3").Copy() | Out-Null
$summary.Range("A2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

for ($r = 3; $r -le 8; $r++) {
    $summary.Range("A$r").Value = $r - 2
}

# ---------------------------------------------------------------------
# 3) Restore the original active-sheet state. Adding "2022-Q3" made it the
#    active sheet; put the selection back on the last sheet ("2021-Q1"),
#    which was the selected tab before this edit.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
